$wb = $excel.ActiveWorkbook

$newNames = @(
    "summ29467270",
    "summ30189971",
    "summ30893713",
    "summ31592945",
    "summ32322825",
    "summ33008604",
    "summ33733539",
    "summ34453102",
    "summ35163028"
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}
